# Corrigindo a leitura e tratamento da tabela e valores vazios
# Update specific "100" placeholder values with the corrected computed values
# for columns G (E), K (I) and N (L) across rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 99.52
$ws.Range("K2").Value = 0.59
$ws.Range("N2").Value = 0.91

$ws.Range("G3").Value = 99.05
$ws.Range("K3").Value = 9.84
$ws.Range("N3").Value = 5.4

$ws.Range("G4").Value = 99.70999999999999
$ws.Range("K4").Value = 12.88
$ws.Range("N4").Value = 6.22

$ws.Range("K5").Value = 66.18000000000001
$ws.Range("N5").Value = 63.27
